$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BitácoraExperiencia1")

# Fill in the two new activity log rows (item 7 and item 8)
$ws.Range("C30").Value = "Se mejora problemas en el html galeria"
$ws.Range("D30").Value = "Nicolas Venegas"

$ws.Range("C31").Value = "Se implementa boton de hamburguesa en los html"
$ws.Range("D31").Value = "Nicolas Venegas"

# Update the active selection to match the saved cursor position
$ws.Activate()
$ws.Range("C32").Select()
